# DPLKLib_Report.xlsx edit:
#  - Remove the PASSWORD_ICONS / KODE_CABANG / NOMOR_TERMINAL columns
#    (F:H) entirely, shifting the remaining columns left.
#  - Update the URL_ICONS hyperlink cell (D2) text from
#    http://192.168.168.107/ to http://192.168.168.111/.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three whole columns F, G, H (PASSWORD_ICONS, KODE_CABANG,
# NOMOR_TERMINAL) - everything to the right shifts left by 3 columns.
$ws.Range("F1:H1").EntireColumn.Delete()

# Update the displayed/stored text of the D2 hyperlink cell.
$ws.Range("D2").Value = "http://192.168.168.111/"
